# Apply the edit: expand the single introductory paragraph into the full
# "Algoritmo del procedimiento" procedure description (multiple paragraphs),
# matching the target OOXML produced by the author's revision.

$d = $word.ActiveDocument

# The whole new body content (12 paragraphs) that replaces the original,
# single, paragraph. Using Range.InsertXML lets us specify the exact
# run/proofErr/language markup Word itself would have produced.
$newBodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Algoritmo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>procedimiento</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ingresar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nombre</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Verificar que el nombre </w:t></w:r><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">no sea RESET sino </w:t></w:r><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>culminar programa</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Ingresar edad </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>Si edad menor a 18 se terminara el programa</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Si edad es mayor a 18 se </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>continuara</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve"> ejecutando</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Si edad esta entre 18 y </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>24  se</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve"> generara un recargo </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Si edad este entre </w:t></w:r><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>25 a 49 se generara un recargo mayor al anterior</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Si edad es mayor o igual a 50 se </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>generara</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve"> un recargo aun mayor a los anteri</w:t></w:r><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>ores</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>Evaluar si tiene esposa</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve">Si tiene esposa se </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>generara</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t xml:space="preserve"> el programa similar a</w:t></w:r><w:r><w:rPr><w:lang w:val="es-GT"/></w:rPr><w:t>l anterior procedimiento</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-GT"/></w:rPr></w:pPr></w:p>
'@

$firstPara = $d.Paragraphs.Item(1)
$rng = $firstPara.Range
$rng.InsertXML($newBodyXml)
